# Edit script generated to apply the diff to poland_division-2_2023-2024.xlsx
# 1) Four row pairs had their match-data columns (F:V) swapped (rows were
#    reordered upstream; Indice/pais/torneio/temporada/data_partida in A:E
#    stay put per row position).
# 2) Rows 91/92/93 had a 3-way rotation of their F:V content.
# 3) A new match row (172) was appended at the end of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap F:V content between row 26 and row 27
$ws.Range("F26").Value = "Hutnik Krakow"
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = "Zaglebie II"
$ws.Range("I26").Value = 1
$ws.Range("J26").Value = 1.89
$ws.Range("K26").Value = "05/08/2023 22:29"
$ws.Range("L26").Value = 1.94
$ws.Range("M26").Value = "06/08/2023 16:50"
$ws.Range("N26").Value = 3.66
$ws.Range("O26").Value = "05/08/2023 22:29"
$ws.Range("P26").Value = 3.54
$ws.Range("Q26").Value = "06/08/2023 16:50"
$ws.Range("R26").Value = 3.62
$ws.Range("S26").Value = "05/08/2023 22:29"
$ws.Range("T26").Value = 3.65
$ws.Range("U26").Value = "06/08/2023 16:50"
$ws.Range("V26").Value = "https://www.betexplorer.com/football/poland/division-2/hutnik-krakow-zaglebie/xdUUh9KF/"
$ws.Range("F27").Value = "Stezyca"
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = "Sandecja Nowy S."
$ws.Range("I27").Value = 1
$ws.Range("J27").Value = 2.27
$ws.Range("K27").Value = "05/08/2023 22:29"
$ws.Range("L27").Value = 2.32
$ws.Range("M27").Value = "06/08/2023 16:44"
$ws.Range("N27").Value = 3.23
$ws.Range("O27").Value = "05/08/2023 22:29"
$ws.Range("P27").Value = 3.31
$ws.Range("Q27").Value = "06/08/2023 15:54"
$ws.Range("R27").Value = 3.05
$ws.Range("S27").Value = "05/08/2023 22:29"
$ws.Range("T27").Value = 2.96
$ws.Range("U27").Value = "06/08/2023 16:44"
$ws.Range("V27").Value = "https://www.betexplorer.com/football/poland/division-2/stezyca-sandecja-nowy-s/WITYiTZL/"

# Swap F:V content between row 41 and row 42
$ws.Range("F41").Value = "Kotwica Kolobrzeg"
$ws.Range("G41").Value = 1
$ws.Range("H41").Value = "Ol. Grudziadz"
$ws.Range("I41").Value = 2
$ws.Range("J41").Value = 1.91
$ws.Range("K41").Value = "19/08/2023 08:43"
$ws.Range("L41").Value = 2.08
$ws.Range("M41").Value = "19/08/2023 17:46"
$ws.Range("N41").Value = 3.45
$ws.Range("O41").Value = "19/08/2023 08:43"
$ws.Range("P41").Value = 3.27
$ws.Range("Q41").Value = "19/08/2023 17:46"
$ws.Range("R41").Value = 3.77
$ws.Range("S41").Value = "19/08/2023 08:43"
$ws.Range("T41").Value = 3.51
$ws.Range("U41").Value = "19/08/2023 17:46"
$ws.Range("V41").Value = "https://www.betexplorer.com/football/poland/division-2/kotwica-kolobrzeg-ol-grudziadz/KIUHC3Bl/"
$ws.Range("F42").Value = "Olimpia Elblag"
$ws.Range("G42").Value = 0
$ws.Range("H42").Value = "Stezyca"
$ws.Range("I42").Value = 1
$ws.Range("J42").Value = 2.27
$ws.Range("K42").Value = "19/08/2023 08:43"
$ws.Range("L42").Value = 2.27
$ws.Range("M42").Value = "19/08/2023 08:43"
$ws.Range("N42").Value = 3.18
$ws.Range("O42").Value = "19/08/2023 08:43"
$ws.Range("P42").Value = 3.2
$ws.Range("Q42").Value = "19/08/2023 16:05"
$ws.Range("R42").Value = 3.09
$ws.Range("S42").Value = "19/08/2023 08:43"
$ws.Range("T42").Value = 3.09
$ws.Range("U42").Value = "19/08/2023 08:43"
$ws.Range("V42").Value = "https://www.betexplorer.com/football/poland/division-2/olimpia-elblag-stezyca/OYTnIPmR/"

# Swap F:V content between row 51 and row 52
$ws.Range("F51").Value = "Pogon Siedlce"
$ws.Range("G51").Value = 3
$ws.Range("H51").Value = "Sandecja Nowy S."
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 2.28
$ws.Range("K51").Value = "26/08/2023 13:13"
$ws.Range("L51").Value = 2.77
$ws.Range("M51").Value = "26/08/2023 16:56"
$ws.Range("N51").Value = 3.08
$ws.Range("O51").Value = "26/08/2023 13:13"
$ws.Range("P51").Value = 3.12
$ws.Range("Q51").Value = "26/08/2023 16:55"
$ws.Range("R51").Value = 3.09
$ws.Range("S51").Value = "26/08/2023 13:13"
$ws.Range("T51").Value = 2.56
$ws.Range("U51").Value = "26/08/2023 16:56"
$ws.Range("V51").Value = "https://www.betexplorer.com/football/poland/division-2/pogon-siedlce-sandecja-nowy-s/bT3PVrIQ/"
$ws.Range("F52").Value = "GKS Jastrzebie"
$ws.Range("G52").Value = 4
$ws.Range("H52").Value = "Polonia Bytom"
$ws.Range("I52").Value = 2
$ws.Range("J52").Value = 1.85
$ws.Range("K52").Value = "26/08/2023 13:13"
$ws.Range("L52").Value = 1.83
$ws.Range("M52").Value = "26/08/2023 16:58"
$ws.Range("N52").Value = 3.42
$ws.Range("O52").Value = "26/08/2023 13:13"
$ws.Range("P52").Value = 3.6
$ws.Range("Q52").Value = "26/08/2023 16:58"
$ws.Range("R52").Value = 4.09
$ws.Range("S52").Value = "26/08/2023 13:13"
$ws.Range("T52").Value = 4.03
$ws.Range("U52").Value = "26/08/2023 16:58"
$ws.Range("V52").Value = "https://www.betexplorer.com/football/poland/division-2/gks-jastrzebie-polonia-bytom/464LWO2K/"

# Swap F:V content between row 148 and row 150
$ws.Range("F148").Value = "Lech Poznan II"
$ws.Range("G148").Value = 1
$ws.Range("H148").Value = "Polonia Bytom"
$ws.Range("I148").Value = 0
$ws.Range("J148").Value = 2.89
$ws.Range("K148").Value = "11/11/2023 01:13"
$ws.Range("L148").Value = 2.98
$ws.Range("M148").Value = "12/11/2023 12:51"
$ws.Range("N148").Value = 3.36
$ws.Range("O148").Value = "11/11/2023 01:13"
$ws.Range("P148").Value = 3.55
$ws.Range("Q148").Value = "12/11/2023 12:51"
$ws.Range("R148").Value = 2.17
$ws.Range("S148").Value = "11/11/2023 01:13"
$ws.Range("T148").Value = 2.21
$ws.Range("U148").Value = "12/11/2023 12:51"
$ws.Range("V148").Value = "https://www.betexplorer.com/football/poland/division-2/lech-poznan-polonia-bytom/8IybC9Ze/"
$ws.Range("F150").Value = "Zaglebie II"
$ws.Range("G150").Value = 4
$ws.Range("H150").Value = "S. Wola"
$ws.Range("I150").Value = 0
$ws.Range("J150").Value = 2.3
$ws.Range("K150").Value = "11/11/2023 01:13"
$ws.Range("L150").Value = 2.72
$ws.Range("M150").Value = "12/11/2023 12:51"
$ws.Range("N150").Value = 3.19
$ws.Range("O150").Value = "11/11/2023 01:13"
$ws.Range("P150").Value = 3.15
$ws.Range("Q150").Value = "12/11/2023 12:51"
$ws.Range("R150").Value = 2.81
$ws.Range("S150").Value = "11/11/2023 01:13"
$ws.Range("T150").Value = 2.59
$ws.Range("U150").Value = "12/11/2023 12:51"
$ws.Range("V150").Value = "https://www.betexplorer.com/football/poland/division-2/zaglebie-stal-stalowa-wola/zqWeXABL/"

# Rotate F:V content among rows 91, 92, 93 (91<-92, 92<-93, 93<-91)
$ws.Range("F91").Value = "Lech Poznan II"
$ws.Range("G91").Value = 0
$ws.Range("H91").Value = "LKS Lodz II"
$ws.Range("I91").Value = 3
$ws.Range("J91").Value = 2.36
$ws.Range("K91").Value = "29/09/2023 02:12"
$ws.Range("L91").Value = 2.4
$ws.Range("M91").Value = "30/09/2023 14:59"
$ws.Range("N91").Value = 3.25
$ws.Range("O91").Value = "29/09/2023 02:12"
$ws.Range("P91").Value = 3.65
$ws.Range("Q91").Value = "30/09/2023 14:59"
$ws.Range("R91").Value = 2.6
$ws.Range("S91").Value = "29/09/2023 02:12"
$ws.Range("T91").Value = 2.63
$ws.Range("U91").Value = "30/09/2023 14:58"
$ws.Range("V91").Value = "https://www.betexplorer.com/football/poland/division-2/lech-poznan-lks-lodz/CdkCuE2k/"
$ws.Range("F92").Value = "Hutnik Krakow"
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = "Chojniczanka"
$ws.Range("I92").Value = 2
$ws.Range("J92").Value = 2.06
$ws.Range("K92").Value = "29/09/2023 02:12"
$ws.Range("L92").Value = 2.24
$ws.Range("M92").Value = "30/09/2023 14:43"
$ws.Range("N92").Value = 3.21
$ws.Range("O92").Value = "29/09/2023 02:12"
$ws.Range("P92").Value = 3.49
$ws.Range("Q92").Value = "30/09/2023 14:41"
$ws.Range("R92").Value = 3.14
$ws.Range("S92").Value = "29/09/2023 02:12"
$ws.Range("T92").Value = 2.95
$ws.Range("U92").Value = "30/09/2023 14:43"
$ws.Range("V92").Value = "https://www.betexplorer.com/football/poland/division-2/hutnik-krakow-chojniczanka/t8BTmxe9/"
$ws.Range("F93").Value = "Sandecja Nowy S."
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = "KKS Kalisz"
$ws.Range("I93").Value = 2
$ws.Range("J93").Value = 2.75
$ws.Range("K93").Value = "29/09/2023 02:12"
$ws.Range("L93").Value = 2.78
$ws.Range("M93").Value = "30/09/2023 14:41"
$ws.Range("N93").Value = 3.09
$ws.Range("O93").Value = "29/09/2023 02:12"
$ws.Range("P93").Value = 3.12
$ws.Range("Q93").Value = "30/09/2023 14:41"
$ws.Range("R93").Value = 2.35
$ws.Range("S93").Value = "29/09/2023 02:12"
$ws.Range("T93").Value = 2.56
$ws.Range("U93").Value = "30/09/2023 14:41"
$ws.Range("V93").Value = "https://www.betexplorer.com/football/poland/division-2/sandecja-nowy-s-kks-kalisz/KbAXndAF/"


# --- Append new row 172 (a new match result) ---
# Copy formatting from the last existing data row (171) down to the new row
# so the index column keeps its bold/centered/bordered style and the date
# column keeps its datetime number format.
$ws.Range("A171:V171").Copy() | Out-Null
$ws.Range("A172:V172").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A172").Value = 171
$ws.Range("B172").Value = "poland"
$ws.Range("C172").Value = "division-2"
$ws.Range("D172").Value = "2023-2024"
$ws.Range("E172").Value = 45263.72916666666
$ws.Range("F172").Value = "Kotwica Kolobrzeg"
$ws.Range("G172").Value = 1
$ws.Range("H172").Value = "Lech Poznan II"
$ws.Range("I172").Value = 2
$ws.Range("J172").Value = 1.41
$ws.Range("K172").Value = "02/12/2023 05:42"
$ws.Range("L172").Value = 1.44
$ws.Range("M172").Value = "03/12/2023 17:19"
$ws.Range("N172").Value = 4.55
$ws.Range("O172").Value = "02/12/2023 05:42"
$ws.Range("P172").Value = 4.89
$ws.Range("Q172").Value = "03/12/2023 17:19"
$ws.Range("R172").Value = 5.7
$ws.Range("S172").Value = "02/12/2023 05:42"
$ws.Range("T172").Value = 5.76
$ws.Range("U172").Value = "03/12/2023 17:19"
$ws.Range("V172").Value = "https://www.betexplorer.com/football/poland/division-2/kotwica-kolobrzeg-lech-poznan/YyXTnBS5/"

Write-Output "Edit complete"
